$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# A new handoff report was generated: the "0bb40168-..." localization file
# moved from "Handed back: in sync with en-US" to "Ready for handoff" (with a
# fresh handoff timestamp), while the "8a655b4d-..." file's entry -- which
# used to be listed second -- is now listed first (still "Handed back: in
# sync with en-US"). Update cell values and hyperlink display text (the
# hyperlink targets / r:id's themselves are untouched) on all three sheets.
# ---------------------------------------------------------------------------

# ===========================================================================
# Overview sheet
# ===========================================================================
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "8a655b4d-52ca-4cc7-af93-86b65082ca1e.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "2016-03-22 20:50:04"

$ws.Range("A3").Value = "0bb40168-d008-4845-bb69-20d061646237.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-03-22 20:51:45"

$map = @{
    '$A$2' = '8a655b4d-52ca-4cc7-af93-86b65082ca1e.md'
    '$A$3' = '0bb40168-d008-4845-bb69-20d061646237.md'
}
foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($map.ContainsKey($addr)) {
        $hl.TextToDisplay = $map[$addr]
    }
}

# ===========================================================================
# zh-cn sheet
# ===========================================================================
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "8a655b4d-52ca-4cc7-af93-86b65082ca1e.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "8a655b4d-52ca-4cc7-af93-86b65082ca1e.8aea50e4078855cf569182d0dd4e83b1b54c45e5.zh-cn.xlf"
$ws.Range("E2").Value = "2016-03-22 20:49:58"
$ws.Range("F2").Value = "8a655b4d-52ca-4cc7-af93-86b65082ca1e.md"
$ws.Range("G2").Value = "8a655b4d-52ca-4cc7-af93-86b65082ca1e.8aea50e4078855cf569182d0dd4e83b1b54c45e5.zh-cn.xlf"
$ws.Range("H2").Value = "2016-03-22 20:50:42"
$ws.Range("J2").Value = "Include"

$ws.Range("A3").Value = "0bb40168-d008-4845-bb69-20d061646237.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "0bb40168-d008-4845-bb69-20d061646237.e99c22eeb7e31ff1578b984d02edcae015cfb77c.zh-cn.xlf"
$ws.Range("E3").Value = "2016-03-22 20:51:41"
$ws.Range("F3").Value = "0bb40168-d008-4845-bb69-20d061646237.md"
$ws.Range("G3").Value = "0bb40168-d008-4845-bb69-20d061646237.e99c22eeb7e31ff1578b984d02edcae015cfb77c.zh-cn.xlf"
$ws.Range("H3").Value = "2016-03-22 20:50:42"
$ws.Range("J3").Value = "Include"

$map = @{
    '$A$2' = '8a655b4d-52ca-4cc7-af93-86b65082ca1e.md'
    '$D$2' = '8a655b4d-52ca-4cc7-af93-86b65082ca1e.8aea50e4078855cf569182d0dd4e83b1b54c45e5.zh-cn.xlf'
    '$F$2' = '8a655b4d-52ca-4cc7-af93-86b65082ca1e.md'
    '$G$2' = '8a655b4d-52ca-4cc7-af93-86b65082ca1e.8aea50e4078855cf569182d0dd4e83b1b54c45e5.zh-cn.xlf'
    '$A$3' = '0bb40168-d008-4845-bb69-20d061646237.md'
    '$D$3' = '0bb40168-d008-4845-bb69-20d061646237.e99c22eeb7e31ff1578b984d02edcae015cfb77c.zh-cn.xlf'
    '$F$3' = '0bb40168-d008-4845-bb69-20d061646237.md'
    '$G$3' = '0bb40168-d008-4845-bb69-20d061646237.e99c22eeb7e31ff1578b984d02edcae015cfb77c.zh-cn.xlf'
}
foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($map.ContainsKey($addr)) {
        $hl.TextToDisplay = $map[$addr]
    }
}

# ===========================================================================
# de-de sheet
# ===========================================================================
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "8a655b4d-52ca-4cc7-af93-86b65082ca1e.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "8a655b4d-52ca-4cc7-af93-86b65082ca1e.8aea50e4078855cf569182d0dd4e83b1b54c45e5.de-de.xlf"
$ws.Range("E2").Value = "2016-03-22 20:50:04"
$ws.Range("F2").Value = "8a655b4d-52ca-4cc7-af93-86b65082ca1e.md"
$ws.Range("G2").Value = "8a655b4d-52ca-4cc7-af93-86b65082ca1e.8aea50e4078855cf569182d0dd4e83b1b54c45e5.de-de.xlf"
$ws.Range("H2").Value = "2016-03-22 20:50:54"
$ws.Range("J2").Value = "Include"

$ws.Range("A3").Value = "0bb40168-d008-4845-bb69-20d061646237.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "0bb40168-d008-4845-bb69-20d061646237.e99c22eeb7e31ff1578b984d02edcae015cfb77c.de-de.xlf"
$ws.Range("E3").Value = "2016-03-22 20:51:45"
$ws.Range("F3").Value = "0bb40168-d008-4845-bb69-20d061646237.md"
$ws.Range("G3").Value = "0bb40168-d008-4845-bb69-20d061646237.e99c22eeb7e31ff1578b984d02edcae015cfb77c.de-de.xlf"
$ws.Range("H3").Value = "2016-03-22 20:50:54"
$ws.Range("J3").Value = "Include"

$map = @{
    '$A$2' = '8a655b4d-52ca-4cc7-af93-86b65082ca1e.md'
    '$D$2' = '8a655b4d-52ca-4cc7-af93-86b65082ca1e.8aea50e4078855cf569182d0dd4e83b1b54c45e5.de-de.xlf'
    '$F$2' = '8a655b4d-52ca-4cc7-af93-86b65082ca1e.md'
    '$G$2' = '8a655b4d-52ca-4cc7-af93-86b65082ca1e.8aea50e4078855cf569182d0dd4e83b1b54c45e5.de-de.xlf'
    '$A$3' = '0bb40168-d008-4845-bb69-20d061646237.md'
    '$D$3' = '0bb40168-d008-4845-bb69-20d061646237.e99c22eeb7e31ff1578b984d02edcae015cfb77c.de-de.xlf'
    '$F$3' = '0bb40168-d008-4845-bb69-20d061646237.md'
    '$G$3' = '0bb40168-d008-4845-bb69-20d061646237.e99c22eeb7e31ff1578b984d02edcae015cfb77c.de-de.xlf'
}
foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($map.ContainsKey($addr)) {
        $hl.TextToDisplay = $map[$addr]
    }
}
